# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.702.53'
$ws.Range('E2').Value = '  +1.83%  '
# Row 3
$ws.Range('D3').Value = '1.897.60'
$ws.Range('E3').Value = '  +2.37%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.19%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.53'
$ws.Range('E5').Value = '  +1.76%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.14%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4811'
$ws.Range('E7').Value = '  +1.05%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2841'
$ws.Range('E8').Value = '  +1.11%  '
# Row 9
$ws.Range('E9').Value = '  +0.81%  '
# Row 10
$ws.Range('D10').Value = '1.882.85'
$ws.Range('E10').Value = '  +1.40%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07485'
$ws.Range('E11').Value = '  +1.89%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.66'
$ws.Range('E12').Value = '  +2.21%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.111'
$ws.Range('E13').Value = '  -0.44%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.18'
$ws.Range('E14').Value = '  +1.17%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6689'
$ws.Range('E15').Value = '  +3.97%  '
# Row 16
$ws.Range('D16').Value = '30.682.30'
$ws.Range('E16').Value = '  +1.92%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.33'
$ws.Range('E17').Value = '  +0.72%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.19%  '
# Row 19
$ws.Range('D19').Value = '2.220.92'
$ws.Range('E19').Value = '  +5.45%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007624'
$ws.Range('E20').Value = '  +0.12%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.24'
$ws.Range('E21').Value = '  +5.19%  '
# Row 22
$ws.Range('E22').Value = '  +1.19%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.13%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.223'
$ws.Range('E24').Value = '  +2.05%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.00'
$ws.Range('E25').Value = '  +2.68%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.277'
$ws.Range('E26').Value = '  -0.03%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.68'
$ws.Range('E27').Value = '  +0.91%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.955'
$ws.Range('E28').Value = '  +2.59%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.403'
$ws.Range('E29').Value = '  -1.71%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1009'
$ws.Range('E30').Value = '  +10.47%  '
# Row 31
$ws.Range('E31').Value = '  +2.92%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.031'
$ws.Range('E32').Value = '  +1.58%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05073'
$ws.Range('E33').Value = '  +1.05%  '
# Row 34
$ws.Range('E34').Value = '  +7.66%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7536'
$ws.Range('E35').Value = '  +1.60%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.712'
$ws.Range('E36').Value = '  +0.75%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01887'
$ws.Range('E37').Value = '  +3.81%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.658'
$ws.Range('E38').Value = '  +1.65%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9189'
$ws.Range('E39').Value = '  +2.06%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.084'
$ws.Range('E40').Value = '  +2.27%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.11'
$ws.Range('E41').Value = '  +0.51%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.857'
$ws.Range('E42').Value = '  -1.13%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4299'
$ws.Range('E43').Value = '  +1.37%  '
# Row 44
$ws.Range('E44').Value = '  +0.34%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.426'
$ws.Range('E45').Value = '  -0.11%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.68'
$ws.Range('E46').Value = '  +1.06%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1274'
$ws.Range('E47').Value = '  -2.76%  '
# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.050'
$ws.Range('E48').Value = '  +2.70%  '
# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.487'
$ws.Range('E49').Value = '  -4.51%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.99'
$ws.Range('E50').Value = '  -0.65%  '
# Row 51
$ws.Range('E51').Value = '  -0.40%  '
